$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 115.8
$ws.Cells.Item(2, 10).Value = 50
$ws.Cells.Item(2, 12).Value = 50
$ws.Cells.Item(2, 14).Value = -276
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 13).ClearContents()
$ws.Cells.Item(51, 8).Value = 11675.451
$ws.Cells.Item(51, 9).Value = 17035.625
$ws.Cells.Item(51, 10).Value = 9811.044
$ws.Cells.Item(51, 11).Value = 17035.625
$ws.Cells.Item(51, 12).Value = 9811.044
$ws.Cells.Item(51, 13).Value = -16551.625
$ws.Cells.Item(51, 14).Value = -10779.044
$ws.Cells.Item(62, 8).Value = 100021496
$ws.Cells.Item(62, 9).Value = 166668340
$ws.Cells.Item(62, 10).Value = 51249.5
$ws.Cells.Item(62, 11).Value = 166668340
$ws.Cells.Item(62, 12).Value = 51249.5
$ws.Cells.Item(62, 13).Value = -166667716
$ws.Cells.Item(62, 14).Value = -52497.5
$ws.Cells.Item(65, 8).Value = 100021496
$ws.Cells.Item(65, 9).Value = 166668340
$ws.Cells.Item(65, 10).Value = 51249.5
$ws.Cells.Item(65, 11).Value = 833341700
$ws.Cells.Item(65, 12).Value = 256247.5
$ws.Cells.Item(65, 13).Value = -833338580
$ws.Cells.Item(65, 14).Value = -262487.5
$ws.Cells.Item(70, 8).Value = 12034.077
$ws.Cells.Item(70, 10).Value = 12911.917
$ws.Cells.Item(70, 12).Value = 38735.751
$ws.Cells.Item(70, 14).Value = -39275.751
$ws.Cells.Item(73, 8).Value = 12034.077
$ws.Cells.Item(73, 10).Value = 12911.917
$ws.Cells.Item(73, 12).Value = 38735.751
$ws.Cells.Item(73, 14).Value = -40607.751
$ws.Cells.Item(98, 8).Value = 6730578.5
$ws.Cells.Item(98, 9).Value = 6497048
$ws.Cells.Item(98, 11).Value = 6497048
$ws.Cells.Item(98, 13).Value = -6495550
$ws.Cells.Item(112, 8).Value = 4357858
$ws.Cells.Item(112, 9).Value = 1452
$ws.Cells.Item(112, 10).Value = 4980201.5
$ws.Cells.Item(112, 11).Value = 4356
$ws.Cells.Item(112, 12).Value = 14940604.5
$ws.Cells.Item(112, 13).Value = -3248
$ws.Cells.Item(112, 14).Value = -14942820.5
$ws.Cells.Item(122, 8).Value = 6730578.5
$ws.Cells.Item(122, 9).Value = 6497048
$ws.Cells.Item(122, 11).Value = 19491144
$ws.Cells.Item(122, 13).Value = -19488694
$ws.Cells.Item(125, 8).Value = 933.46155
$ws.Cells.Item(125, 9).Value = 638.75
$ws.Cells.Item(125, 11).Value = 5748.75
$ws.Cells.Item(125, 13).Value = -3288.75
$ws.Cells.Item(129, 8).Value = 1126.7142
$ws.Cells.Item(129, 9).Value = 1126.7142
$ws.Cells.Item(129, 11).Value = 3380.1426
$ws.Cells.Item(129, 13).Value = 1619.8574
$ws.Cells.Item(135, 8).Value = 1678.3334
$ws.Cells.Item(135, 9).Value = 300
$ws.Cells.Item(135, 10).Value = 2367.5
$ws.Cells.Item(135, 11).Value = 2700
$ws.Cells.Item(135, 12).Value = 21307.5
$ws.Cells.Item(135, 13).Value = -165
$ws.Cells.Item(135, 14).Value = -26377.5
$ws.Cells.Item(137, 8).Value = 42729.207
$ws.Cells.Item(137, 9).Value = 59324.05
$ws.Cells.Item(137, 11).Value = 177972.15
$ws.Cells.Item(137, 13).Value = -175422.15
$ws.Cells.Item(138, 8).Value = 2727.83
$ws.Cells.Item(138, 9).Value = 1263.0322
$ws.Cells.Item(138, 10).Value = 3385.9275
$ws.Cells.Item(138, 11).Value = 3789.0966
$ws.Cells.Item(138, 12).Value = 10157.7825
$ws.Cells.Item(138, 13).Value = 1350.9034
$ws.Cells.Item(138, 14).Value = -20437.7825

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9923.1
$ws.Cells.Item(32, 9).Value = 7668.75
$ws.Cells.Item(32, 10).Value = 15720
$ws.Cells.Item(32, 11).Value = 7668.75
$ws.Cells.Item(32, 12).Value = 15720
$ws.Cells.Item(32, 13).Value = -7381.75
$ws.Cells.Item(32, 14).Value = -16294
$ws.Cells.Item(45, 8).Value = 1566.6666
$ws.Cells.Item(45, 9).Value = 1350
$ws.Cells.Item(45, 11).Value = 1350
$ws.Cells.Item(45, 13).Value = -973
$ws.Cells.Item(88, 8).Value = 3217.2856
$ws.Cells.Item(88, 10).Value = 3168.5
$ws.Cells.Item(88, 12).Value = 3168.5
$ws.Cells.Item(88, 14).Value = -3980.5
$ws.Cells.Item(91, 8).Value = 3217.2856
$ws.Cells.Item(91, 10).Value = 3168.5
$ws.Cells.Item(91, 12).Value = 3168.5
$ws.Cells.Item(91, 14).Value = -5976.5
$ws.Cells.Item(97, 8).Value = 28572564
$ws.Cells.Item(97, 9).Value = 33333992
$ws.Cells.Item(97, 10).Value = 4000
$ws.Cells.Item(97, 11).Value = 33333992
$ws.Cells.Item(97, 12).Value = 4000
$ws.Cells.Item(97, 13).Value = -33333496
$ws.Cells.Item(97, 14).Value = -4992
$ws.Cells.Item(110, 8).Value = 13556
$ws.Cells.Item(110, 9).Value = 18994.5
$ws.Cells.Item(110, 10).Value = 5398.25
$ws.Cells.Item(110, 11).Value = 18994.5
$ws.Cells.Item(110, 12).Value = 5398.25
$ws.Cells.Item(110, 13).Value = -16949.5
$ws.Cells.Item(110, 14).Value = -9488.25
$ws.Cells.Item(122, 8).Value = 6890.421
$ws.Cells.Item(122, 9).Value = 4993.4614
$ws.Cells.Item(122, 11).Value = 14980.3842
$ws.Cells.Item(122, 13).Value = -12530.3842
$ws.Cells.Item(132, 8).Value = 1580.5476
$ws.Cells.Item(132, 9).Value = 1625.3243
$ws.Cells.Item(132, 10).Value = 1249.2
$ws.Cells.Item(132, 11).Value = 4875.9729
$ws.Cells.Item(132, 12).Value = 3747.6
$ws.Cells.Item(132, 13).Value = -2345.9729
$ws.Cells.Item(132, 14).Value = -8807.6
$ws.Cells.Item(135, 8).Value = 89999.836
$ws.Cells.Item(135, 10).Value = 89999.836
$ws.Cells.Item(135, 12).Value = 89999.836
$ws.Cells.Item(135, 14).Value = -100139.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 7499.5
$ws.Cells.Item(7, 9).Value = 4999
$ws.Cells.Item(7, 10).Value = 10000
$ws.Cells.Item(7, 11).Value = 4999
$ws.Cells.Item(7, 12).Value = 10000
$ws.Cells.Item(7, 13).Value = -4886
$ws.Cells.Item(7, 14).Value = -10226
$ws.Cells.Item(20, 8).Value = 1217.5714
$ws.Cells.Item(20, 9).Value = 1090.7142
$ws.Cells.Item(20, 10).Value = 1598.1428
$ws.Cells.Item(20, 11).Value = 1090.7142
$ws.Cells.Item(20, 12).Value = 1598.1428
$ws.Cells.Item(20, 13).Value = -843.7141999999999
$ws.Cells.Item(20, 14).Value = -2092.1428
$ws.Cells.Item(26, 8).Value = 210200.8
$ws.Cells.Item(26, 9).Value = 210200.8
$ws.Cells.Item(26, 11).Value = 210200.8
$ws.Cells.Item(26, 13).Value = -209908.8
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 2159.2144
$ws.Cells.Item(86, 9).Value = 2208.6365
$ws.Cells.Item(86, 10).Value = 1978
$ws.Cells.Item(86, 11).Value = 2208.6365
$ws.Cells.Item(86, 12).Value = 1978
$ws.Cells.Item(86, 13).Value = -1085.6365
$ws.Cells.Item(86, 14).Value = -4224
$ws.Cells.Item(89, 8).Value = 2159.2144
$ws.Cells.Item(89, 9).Value = 2208.6365
$ws.Cells.Item(89, 10).Value = 1978
$ws.Cells.Item(89, 11).Value = 11043.1825
$ws.Cells.Item(89, 12).Value = 9890
$ws.Cells.Item(89, 13).Value = -5427.182500000001
$ws.Cells.Item(89, 14).Value = -21122
$ws.Cells.Item(94, 8).Value = 1293.2858
$ws.Cells.Item(94, 9).Value = 1293.2858
$ws.Cells.Item(94, 11).Value = 1293.2858
$ws.Cells.Item(94, 13).Value = -842.2858000000001
$ws.Cells.Item(105, 8).Value = 1762.5714
$ws.Cells.Item(105, 9).Value = 1762.5714
$ws.Cells.Item(105, 11).Value = 1762.5714
$ws.Cells.Item(105, 13).Value = -15.57140000000004

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 419076.25
$ws.Cells.Item(31, 9).Value = 770525.9
$ws.Cells.Item(31, 10).Value = 3726.7273
$ws.Cells.Item(31, 11).Value = 770525.9
$ws.Cells.Item(31, 12).Value = 3726.7273
$ws.Cells.Item(31, 13).Value = -770230.9
$ws.Cells.Item(31, 14).Value = -4316.7273
$ws.Cells.Item(34, 8).Value = 419076.25
$ws.Cells.Item(34, 9).Value = 770525.9
$ws.Cells.Item(34, 10).Value = 3726.7273
$ws.Cells.Item(34, 11).Value = 770525.9
$ws.Cells.Item(34, 12).Value = 3726.7273
$ws.Cells.Item(34, 13).Value = -770323.9
$ws.Cells.Item(34, 14).Value = -4130.7273
$ws.Cells.Item(39, 8).Value = 3200.1428
$ws.Cells.Item(39, 9).Value = 3200.1428
$ws.Cells.Item(39, 11).Value = 3200.1428
$ws.Cells.Item(39, 13).Value = -2809.1428
$ws.Cells.Item(42, 8).Value = 16189.4
$ws.Cells.Item(42, 9).Value = 8450
$ws.Cells.Item(42, 10).Value = 18124.25
$ws.Cells.Item(42, 11).Value = 8450
$ws.Cells.Item(42, 12).Value = 18124.25
$ws.Cells.Item(42, 13).Value = -7857
$ws.Cells.Item(42, 14).Value = -19310.25
$ws.Cells.Item(49, 8).Value = 3200.1428
$ws.Cells.Item(49, 9).Value = 3200.1428
$ws.Cells.Item(49, 11).Value = 3200.1428
$ws.Cells.Item(49, 13).Value = -3018.1428
$ws.Cells.Item(58, 8).Value = 1394.5
$ws.Cells.Item(58, 9).Value = 1438.7727
$ws.Cells.Item(58, 10).Value = 1232.1666
$ws.Cells.Item(58, 11).Value = 1438.7727
$ws.Cells.Item(58, 12).Value = 1232.1666
$ws.Cells.Item(58, 13).Value = -1235.7727
$ws.Cells.Item(58, 14).Value = -1638.1666
$ws.Cells.Item(86, 8).Value = 4225.6
$ws.Cells.Item(86, 9).Value = 4495.8335
$ws.Cells.Item(86, 11).Value = 4495.8335
$ws.Cells.Item(86, 13).Value = -3372.8335
$ws.Cells.Item(89, 8).Value = 4225.6
$ws.Cells.Item(89, 9).Value = 4495.8335
$ws.Cells.Item(89, 11).Value = 22479.1675
$ws.Cells.Item(89, 13).Value = -16863.1675
$ws.Cells.Item(105, 8).Value = 4297.413
$ws.Cells.Item(105, 9).Value = 1828.1305
$ws.Cells.Item(105, 11).Value = 1828.1305
$ws.Cells.Item(105, 13).Value = -81.13049999999998
$ws.Cells.Item(106, 8).Value = 28198.75
$ws.Cells.Item(106, 10).Value = 28198.75
$ws.Cells.Item(106, 12).Value = 28198.75
$ws.Cells.Item(106, 14).Value = -30722.75
$ws.Cells.Item(122, 8).Value = 4589.3
$ws.Cells.Item(122, 9).Value = 2474.5
$ws.Cells.Item(122, 10).Value = 5999.1665
$ws.Cells.Item(122, 11).Value = 7423.5
$ws.Cells.Item(122, 12).Value = 17997.4995
$ws.Cells.Item(122, 13).Value = -4973.5
$ws.Cells.Item(122, 14).Value = -22897.4995
$ws.Cells.Item(136, 8).Value = 1394.5
$ws.Cells.Item(136, 9).Value = 1438.7727
$ws.Cells.Item(136, 10).Value = 1232.1666
$ws.Cells.Item(136, 11).Value = 4316.3181
$ws.Cells.Item(136, 12).Value = 3696.4998
$ws.Cells.Item(136, 13).Value = -1766.3181
$ws.Cells.Item(136, 14).Value = -8796.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 599.4
$ws.Cells.Item(92, 10).Value = 250
$ws.Cells.Item(92, 12).Value = 750
$ws.Cells.Item(92, 14).Value = -3246
$ws.Cells.Item(97, 8).Value = 369.75
$ws.Cells.Item(97, 9).Value = 388.5
$ws.Cells.Item(97, 10).Value = 363.5
$ws.Cells.Item(97, 11).Value = 1165.5
$ws.Cells.Item(97, 12).Value = 1090.5
$ws.Cells.Item(97, 13).Value = -669.5
$ws.Cells.Item(97, 14).Value = -2082.5
$ws.Cells.Item(109, 8).Value = 1400.5
$ws.Cells.Item(109, 9).Value = 1000.625
$ws.Cells.Item(109, 11).Value = 3001.875
$ws.Cells.Item(109, 13).Value = -1961.875
$ws.Cells.Item(114, 8).Value = 639.75
$ws.Cells.Item(114, 9).Value = 639.75
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 11).Value = 1919.25
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 13).Value = 1334.75
$ws.Cells.Item(114, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 4282.5713
$ws.Cells.Item(122, 9).Value = 4497.75
$ws.Cells.Item(122, 10).Value = 3995.6667
$ws.Cells.Item(122, 11).Value = 40479.75
$ws.Cells.Item(122, 12).Value = 35961.0003
$ws.Cells.Item(122, 13).Value = -38029.75
$ws.Cells.Item(122, 14).Value = -40861.0003
$ws.Cells.Item(127, 8).Value = 2039.4
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 2039.4
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 6118.200000000001
$ws.Cells.Item(127, 13).ClearContents()
$ws.Cells.Item(127, 14).Value = -16038.2
$ws.Cells.Item(131, 8).Value = 1332.1428
$ws.Cells.Item(131, 9).Value = 496.1111
$ws.Cells.Item(131, 11).Value = 1488.3333
$ws.Cells.Item(131, 13).Value = 3551.6667
$ws.Cells.Item(133, 8).Value = 12136.667
$ws.Cells.Item(133, 9).Value = 4455
$ws.Cells.Item(133, 10).Value = 27500
$ws.Cells.Item(133, 11).Value = 13365
$ws.Cells.Item(133, 12).Value = 82500
$ws.Cells.Item(133, 13).Value = -8305
$ws.Cells.Item(133, 14).Value = -92620
$ws.Cells.Item(134, 8).Value = 12944.238
$ws.Cells.Item(134, 9).Value = 14695.823
$ws.Cells.Item(134, 10).Value = 5500
$ws.Cells.Item(134, 11).Value = 44087.469
$ws.Cells.Item(134, 12).Value = 16500
$ws.Cells.Item(134, 13).Value = -39017.469
$ws.Cells.Item(134, 14).Value = -26640
$ws.Cells.Item(139, 8).Value = 4017.9
$ws.Cells.Item(139, 9).Value = 3647.375
$ws.Cells.Item(139, 11).Value = 10942.125
$ws.Cells.Item(139, 13).Value = -5802.125
$ws.Cells.Item(141, 8).Value = 3695.5
$ws.Cells.Item(141, 9).Value = 3094
$ws.Cells.Item(141, 11).Value = 9282
$ws.Cells.Item(141, 13).Value = -4102

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 43.81818
$ws.Cells.Item(2, 9).Value = 44.75
$ws.Cells.Item(2, 10).Value = 41.333332
$ws.Cells.Item(2, 11).Value = 44.75
$ws.Cells.Item(2, 12).Value = 41.333332
$ws.Cells.Item(2, 13).Value = 68.25
$ws.Cells.Item(2, 14).Value = -267.333332
$ws.Cells.Item(15, 8).Value = 13814.75
$ws.Cells.Item(15, 10).Value = 8420
$ws.Cells.Item(15, 12).Value = 8420
$ws.Cells.Item(15, 14).Value = -8996
$ws.Cells.Item(35, 8).Value = 14166.667
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 13).ClearContents()
$ws.Cells.Item(70, 8).Value = 9874.25
$ws.Cells.Item(70, 9).Value = 9862.817999999999
$ws.Cells.Item(70, 11).Value = 9862.817999999999
$ws.Cells.Item(70, 13).Value = -9592.817999999999
$ws.Cells.Item(73, 8).Value = 9874.25
$ws.Cells.Item(73, 9).Value = 9862.817999999999
$ws.Cells.Item(73, 11).Value = 9862.817999999999
$ws.Cells.Item(73, 13).Value = -8926.817999999999
$ws.Cells.Item(80, 8).Value = 4559.074
$ws.Cells.Item(80, 9).Value = 3996
$ws.Cells.Item(80, 11).Value = 3996
$ws.Cells.Item(80, 13).Value = -2998
$ws.Cells.Item(81, 8).Value = 13814.75
$ws.Cells.Item(81, 10).Value = 8420
$ws.Cells.Item(81, 12).Value = 8420
$ws.Cells.Item(81, 14).Value = -10416
$ws.Cells.Item(83, 8).Value = 4559.074
$ws.Cells.Item(83, 9).Value = 3996
$ws.Cells.Item(83, 11).Value = 19980
$ws.Cells.Item(83, 13).Value = -14988
$ws.Cells.Item(84, 8).Value = 13814.75
$ws.Cells.Item(84, 10).Value = 8420
$ws.Cells.Item(84, 12).Value = 25260
$ws.Cells.Item(84, 14).Value = -35244
$ws.Cells.Item(97, 8).Value = 5251.625
$ws.Cells.Item(97, 9).Value = 3829.75
$ws.Cells.Item(97, 10).Value = 6673.5
$ws.Cells.Item(97, 11).Value = 3829.75
$ws.Cells.Item(97, 12).Value = 6673.5
$ws.Cells.Item(97, 13).Value = -3333.75
$ws.Cells.Item(97, 14).Value = -7665.5
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 3634.6562
$ws.Cells.Item(122, 9).Value = 3569.44
$ws.Cells.Item(122, 10).Value = 3867.5715
$ws.Cells.Item(122, 11).Value = 10708.32
$ws.Cells.Item(122, 12).Value = 11602.7145
$ws.Cells.Item(122, 13).Value = -8258.32
$ws.Cells.Item(122, 14).Value = -16502.7145
$ws.Cells.Item(123, 8).Value = 31208.5
$ws.Cells.Item(123, 10).Value = 31208.5
$ws.Cells.Item(123, 12).Value = 31208.5
$ws.Cells.Item(123, 14).Value = -36108.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 6004.7144
$ws.Cells.Item(40, 9).Value = 7842
$ws.Cells.Item(40, 11).Value = 7842
$ws.Cells.Item(40, 13).Value = -7706
$ws.Cells.Item(41, 8).Value = 18999.334
$ws.Cells.Item(41, 9).Value = 18999.334
$ws.Cells.Item(41, 11).Value = 18999.334
$ws.Cells.Item(41, 13).Value = -18561.334
$ws.Cells.Item(46, 8).Value = 2912.5908
$ws.Cells.Item(46, 9).Value = 2082.8462
$ws.Cells.Item(46, 10).Value = 4111.1113
$ws.Cells.Item(46, 11).Value = 2082.8462
$ws.Cells.Item(46, 12).Value = 4111.1113
$ws.Cells.Item(46, 13).Value = -1894.8462
$ws.Cells.Item(46, 14).Value = -4487.1113
$ws.Cells.Item(47, 8).Value = 35500
$ws.Cells.Item(47, 10).Value = 35500
$ws.Cells.Item(47, 12).Value = 35500
$ws.Cells.Item(47, 14).Value = -36480
$ws.Cells.Item(52, 8).Value = 35500
$ws.Cells.Item(52, 10).Value = 35500
$ws.Cells.Item(52, 12).Value = 35500
$ws.Cells.Item(52, 14).Value = -35966
$ws.Cells.Item(61, 8).Value = 2684
$ws.Cells.Item(61, 9).Value = 2694.6365
$ws.Cells.Item(61, 10).Value = 2645
$ws.Cells.Item(61, 11).Value = 2694.6365
$ws.Cells.Item(61, 12).Value = 2645
$ws.Cells.Item(61, 13).Value = -2492.6365
$ws.Cells.Item(61, 14).Value = -3049
$ws.Cells.Item(68, 8).Value = 8874.833000000001
$ws.Cells.Item(68, 9).Value = 10684
$ws.Cells.Item(68, 10).Value = 2000
$ws.Cells.Item(68, 11).Value = 10684
$ws.Cells.Item(68, 12).Value = 2000
$ws.Cells.Item(68, 13).Value = -9935
$ws.Cells.Item(68, 14).Value = -3498
$ws.Cells.Item(71, 8).Value = 8874.833000000001
$ws.Cells.Item(71, 9).Value = 10684
$ws.Cells.Item(71, 10).Value = 2000
$ws.Cells.Item(71, 11).Value = 53420
$ws.Cells.Item(71, 12).Value = 10000
$ws.Cells.Item(71, 13).Value = -49676
$ws.Cells.Item(71, 14).Value = -17488
$ws.Cells.Item(82, 8).Value = 2532.875
$ws.Cells.Item(82, 9).Value = 2480
$ws.Cells.Item(82, 11).Value = 2480
$ws.Cells.Item(82, 13).Value = -2119
$ws.Cells.Item(85, 8).Value = 2532.875
$ws.Cells.Item(85, 9).Value = 2480
$ws.Cells.Item(85, 11).Value = 2480
$ws.Cells.Item(85, 13).Value = -1232
$ws.Cells.Item(93, 8).Value = 66670680
$ws.Cells.Item(93, 10).Value = 125004620
$ws.Cells.Item(93, 12).Value = 125004620
$ws.Cells.Item(93, 14).Value = -125007116
$ws.Cells.Item(113, 8).Value = 2684
$ws.Cells.Item(113, 9).Value = 2694.6365
$ws.Cells.Item(113, 10).Value = 2645
$ws.Cells.Item(113, 11).Value = 2694.6365
$ws.Cells.Item(113, 12).Value = 2645
$ws.Cells.Item(113, 13).Value = -524.6365000000001
$ws.Cells.Item(113, 14).Value = -6985
$ws.Cells.Item(122, 8).Value = 5779.8125
$ws.Cells.Item(122, 9).Value = 6079.8184
$ws.Cells.Item(122, 11).Value = 18239.4552
$ws.Cells.Item(122, 13).Value = -15789.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(49, 8).Value = 49999
$ws.Cells.Item(49, 9).Value = 49999
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 11).Value = 49999
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 13).Value = -49769
$ws.Cells.Item(49, 14).ClearContents()
$ws.Cells.Item(51, 8).Value = 19998.8
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 14).ClearContents()
$ws.Cells.Item(81, 8).Value = 624.25
$ws.Cells.Item(81, 9).Value = 499
$ws.Cells.Item(81, 11).Value = 998
$ws.Cells.Item(81, 13).Value = 63
$ws.Cells.Item(84, 8).Value = 624.25
$ws.Cells.Item(84, 9).Value = 499
$ws.Cells.Item(84, 11).Value = 4990
$ws.Cells.Item(84, 13).Value = 314
$ws.Cells.Item(96, 8).Value = 145743.28
$ws.Cells.Item(96, 10).Value = 4100
$ws.Cells.Item(96, 12).Value = 4100
$ws.Cells.Item(96, 14).Value = -6846
$ws.Cells.Item(100, 8).Value = 933.2
$ws.Cells.Item(100, 9).Value = 863.7857
$ws.Cells.Item(100, 10).Value = 1095.1666
$ws.Cells.Item(100, 11).Value = 1727.5714
$ws.Cells.Item(100, 12).Value = 2190.3332
$ws.Cells.Item(100, 13).Value = -1186.5714
$ws.Cells.Item(100, 14).Value = -3272.3332
$ws.Cells.Item(126, 8).Value = 2139.5833
$ws.Cells.Item(126, 9).Value = 2088.818
$ws.Cells.Item(126, 11).Value = 6266.454000000001
$ws.Cells.Item(126, 13).Value = -3796.454000000001
$ws.Cells.Item(132, 8).Value = 1127.1316
$ws.Cells.Item(132, 9).Value = 1127.1316
$ws.Cells.Item(132, 11).Value = 3381.3948
$ws.Cells.Item(132, 13).Value = -851.3948
